# Weekly price update for Fruta / Hortaliza (Haba) dataset.
# Two new weekly records are inserted into the data table:
#   - one before the current row 279 (pushes the existing 279..357 rows down by one)
#   - one before the (post-shift) row 356, i.e. right before the record that used to be
#     row 355, so the tail of the table (old rows 355..357) ends up at 357..359.
# Both inserted rows replicate the recurring "Provincia de Limarí" price point for the
# relevant quality grade, stamped with the next date(s) in the weekly series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-HabaRow {
    param($RowNumber, $Fecha, $Calidad, $Volumen, $PrecioMinimo, $PrecioMaximo, $PrecioPromedio, $Origen, $PrecioKg)

    $ws.Cells.Item($RowNumber, 1).Value = 9
    $ws.Cells.Item($RowNumber, 2).Value = 'Vega Central Mapocho de Santiago'
    $ws.Cells.Item($RowNumber, 3).Value = 'Metropolitana'
    $ws.Cells.Item($RowNumber, 4).Value = $Fecha
    $ws.Cells.Item($RowNumber, 5).Value = 13
    $ws.Cells.Item($RowNumber, 6).Value = 100112026
    $ws.Cells.Item($RowNumber, 7).Value = 'Haba'
    $ws.Cells.Item($RowNumber, 8).Value = 'Sin especificar'
    $ws.Cells.Item($RowNumber, 9).Value = $Calidad
    $ws.Cells.Item($RowNumber, 10).Value = $Volumen
    $ws.Cells.Item($RowNumber, 11).Value = $PrecioMinimo
    $ws.Cells.Item($RowNumber, 12).Value = $PrecioMaximo
    $ws.Cells.Item($RowNumber, 13).Value = $PrecioPromedio
    $ws.Cells.Item($RowNumber, 14).Value = '$/saco 25 kilos'
    $ws.Cells.Item($RowNumber, 15).Value = $Origen
    $ws.Cells.Item($RowNumber, 16).Value = $PrecioKg
    $ws.Cells.Item($RowNumber, 17).Value = 25
    $ws.Cells.Item($RowNumber, 18).Value = 'Hortaliza'
}

# Insert first new weekly record at row 279 (shifts old rows 279..357 -> 280..358).
$ws.Rows.Item(279).Insert()
Set-HabaRow 279 45120 'Primera' 70 15000 16000 15500 'Provincia de Limarí' 620

# Insert second new weekly record at row 356 (shifts the current rows 356..358 -> 357..359).
$ws.Rows.Item(356).Insert()
Set-HabaRow 356 45121 'Primera' 52 15000 16000 15500 'Provincia de Limarí' 620
